$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Recall/Precision columns are treated as text so that full-precision
# decimal strings (and any embedded newlines) are preserved exactly,
# matching the original data export format.
$ws.Range("C2:D41").NumberFormat = "@"

$ws.Range("C2").Value = "0.5"
$ws.Range("D2").Value = "0.1818181818181818181818181818`n"
$ws.Range("C3").Value = "0.3333333333333333333333333333"
$ws.Range("D3").Value = "0.1818181818181818181818181818`n"
$ws.Range("C4").Value = "1"
$ws.Range("D4").Value = "0.09090909090909090909090909091`n"
$ws.Range("C5").Value = "0.3"
$ws.Range("D5").Value = "0.2727272727272727272727272727`n"
$ws.Range("C6").Value = "0.5"
$ws.Range("D6").Value = "0.2727272727272727272727272727`n"
$ws.Range("C7").Value = "0.3333333333333333333333333333"
$ws.Range("D7").Value = "0.1818181818181818181818181818`n"
$ws.Range("C8").Value = "0.5"
$ws.Range("D8").Value = "0.2727272727272727272727272727`n"
$ws.Range("C9").Value = "1"
$ws.Range("D9").Value = "0.09090909090909090909090909091`n"
$ws.Range("C10").Value = "0.3636363636363636363636363636"
$ws.Range("D10").Value = "0.3636363636363636363636363636`n"
$ws.Range("C11").Value = "0.4285714285714285714285714286"
$ws.Range("D11").Value = "0.2727272727272727272727272727`n"
$ws.Range("C12").Value = "1"
$ws.Range("D12").Value = "0.1818181818181818181818181818`n"
$ws.Range("C13").Value = "0.7142857142857142857142857143"
$ws.Range("D13").Value = "0.4545454545454545454545454545`n"
$ws.Range("C14").Value = "0"
$ws.Range("D14").Value = "0`n"
$ws.Range("C15").Value = "0.2941176470588235294117647059"
$ws.Range("D15").Value = "0.4545454545454545454545454545`n"
$ws.Range("C16").Value = "0"
$ws.Range("D16").Value = "0`n"
$ws.Range("C17").Value = "1"
$ws.Range("D17").Value = "0.1818181818181818181818181818`n"
$ws.Range("C18").Value = "0.8333333333333333333333333333"
$ws.Range("D18").Value = "0.4545454545454545454545454545`n"
$ws.Range("C19").Value = "0"
$ws.Range("D19").Value = "0`n"
$ws.Range("C20").Value = "0.2941176470588235294117647059"
$ws.Range("D20").Value = "0.4545454545454545454545454545`n"
$ws.Range("C21").Value = "0.2"
$ws.Range("D21").Value = "0.1818181818181818181818181818`n"
$ws.Range("C22").Value = "0.4"
$ws.Range("D22").Value = "0.1818181818181818181818181818`n"
$ws.Range("C23").Value = "0.8333333333333333333333333333"
$ws.Range("D23").Value = "0.4545454545454545454545454545`n"
$ws.Range("C24").Value = "0"
$ws.Range("D24").Value = "0`n"
$ws.Range("C25").Value = "0.25"
$ws.Range("D25").Value = "0.3636363636363636363636363636`n"
$ws.Range("C26").Value = "0.1"
$ws.Range("D26").Value = "0.09090909090909090909090909091`n"
$ws.Range("C27").Value = "0.3333333333333333333333333333"
$ws.Range("D27").Value = "0.09090909090909090909090909091`n"
$ws.Range("C28").Value = "0.6666666666666666666666666667"
$ws.Range("D28").Value = "0.3636363636363636363636363636`n"
$ws.Range("C29").Value = "0"
$ws.Range("D29").Value = "0`n"
$ws.Range("C30").Value = "0.3333333333333333333333333333"
$ws.Range("D30").Value = "0.5454545454545454545454545455`n"
$ws.Range("C31").Value = "0.1"
$ws.Range("D31").Value = "0.09090909090909090909090909091`n"
$ws.Range("C32").Value = "0.4"
$ws.Range("D32").Value = "0.1818181818181818181818181818`n"
$ws.Range("C33").Value = "0.7142857142857142857142857143"
$ws.Range("D33").Value = "0.4545454545454545454545454545`n"
$ws.Range("C34").Value = "0"
$ws.Range("D34").Value = "0`n"
$ws.Range("C35").Value = "0.3125"
$ws.Range("D35").Value = "0.4545454545454545454545454545`n"
$ws.Range("C36").Value = "0"
$ws.Range("D36").Value = "0`n"
$ws.Range("C37").Value = "0.5"
$ws.Range("D37").Value = "0.1818181818181818181818181818`n"
$ws.Range("C38").Value = "0.7142857142857142857142857143"
$ws.Range("D38").Value = "0.4545454545454545454545454545`n"
$ws.Range("C39").Value = "0"
$ws.Range("D39").Value = "0`n"
$ws.Range("C40").Value = "0.1764705882352941176470588235"
$ws.Range("D40").Value = "0.2727272727272727272727272727`n"
$ws.Range("C41").Value = "0.1"
$ws.Range("D41").Value = "0.09090909090909090909090909091`n"

# Clear the automatic row-height bump that Excel applies when it sees the
# embedded newline characters, restoring the default (non-custom) row height.
$ws.Rows("1:41").AutoFit()

